$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 41312
$ws.Range("G2").Value = "暂时售罄"
$ws.Range("F3").Value = 207
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 9253
$ws.Range("F6").Value = 193
$ws.Range("F7").Value = 755
$ws.Range("F8").Value = 844
$ws.Range("F9").Value = 689
$ws.Range("F10").Value = 190
$ws.Range("F12").Value = 271
$ws.Range("F13").Value = 844
$ws.Range("F14").Value = 76
$ws.Range("F15").Value = 114
$ws.Range("F16").Value = 680
$ws.Range("F17").Value = 289
$ws.Range("F18").Value = 1317
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 576
$ws.Range("F21").Value = 668
$ws.Range("F22").Value = 443
$ws.Range("F23").Value = 653
$ws.Range("F24").Value = 699
$ws.Range("F26").Value = 39
$ws.Range("F27").Value = 54
$ws.Range("F28").Value = 461
$ws.Range("F29").Value = 494
$ws.Range("F30").Value = 36
$ws.Range("F31").Value = 209
$ws.Range("F32").Value = 905
$ws.Range("F33").Value = 4
$ws.Range("F34").Value = 423
$ws.Range("F35").Value = 74
$ws.Range("F37").Value = 133
$ws.Range("F38").Value = 339
$ws.Range("F39").Value = 1192
$ws.Range("F40").Value = 275
$ws.Range("F42").Value = 1192
$ws.Range("F43").Value = 360
$ws.Range("F45").Value = 4
$ws.Range("F46").Value = 14
$ws.Range("F47").Value = 26
$ws.Range("F48").Value = 36
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1974
$ws.Range("F5").Value = 4432
$ws.Range("F13").Value = 53
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1974
$ws.Range("F3").Value = 480
$ws.Range("F4").Value = 327
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 480
$ws.Range("F4").Value = 41313
$ws.Range("G4").Value = 0
$ws.Range("F7").Value = 207
$ws.Range("F9").Value = 15
$ws.Range("F11").Value = 9253
$ws.Range("F12").Value = 193
$ws.Range("F13").Value = 755
$ws.Range("C14").Value = "广州·LookLook动漫嘉年华2th"
$ws.Range("D14").Value = "展贸东路200号 恒达智慧汽车城"
$ws.Range("E14").Value = "2024.07.27 10:00-07.28 17:30"
$ws.Range("F14").Value = 755
$ws.Range("G14").Value = 19.9
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=87217"
$ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202407/wjvAqamr1720170199991.jpeg"
$ws.Range("C15").Value = "广州·fhana ONE MAN LIVE 巡回演唱会 2024"
$ws.Range("D15").Value = "流花路117号流花展贸中心5号馆 广州大麦66live house"
$ws.Range("E15").Value = "2024.07.27 19:00-07.27 21:30"
$ws.Range("F15").Value = 76
$ws.Range("G15").Value = 380
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=87638"
$ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202406/P9TXBIjT1718746868925.jpeg"
$ws.Range("C16").Value = "广州·凹凸世界八周年 夏日特调主题嘉年华"
$ws.Range("D16").Value = "动漫星城 动漫星城"
$ws.Range("E16").Value = "2024.07.27 00:00-07.29 23:59"
$ws.Range("F16").Value = 327
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=88982"
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202407/ED5bP47A1720579350426.jpeg"
$ws.Range("C17").Value = "广州·原神x星穹x崩only"
$ws.Range("D17").Value = "鸿盛二路巨大创意产业园 巨大产业园·智汇港"
$ws.Range("E17").Value = "2024.07.27 10:00-07.27 17:00"
$ws.Range("F17").Value = 844
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=87184"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202406/u67hjpFi1718160712051.jpeg"
$ws.Range("B18").Value = "2024-07-27"
$ws.Range("C18").Value = "广州·广友·星声代起-东方同人嘉年华"
$ws.Range("D18").Value = "黄埔大道中309-315号 羊城创意产业园"
$ws.Range("E18").Value = "2024.07.27 09:00-07.27 22:00"
$ws.Range("F18").Value = 116
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=88303"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202406/4oMyFoF11719567403862.jpeg"
$ws.Range("C19").Value = "广州·小马宝莉only"
$ws.Range("D19").Value = "鸿盛二路巨大创意产业园 巨大产业园·智汇港"
$ws.Range("F19").Value = 271
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=88110"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202406/XH2hCwrg1719393458914.jpeg"
$ws.Range("B20").Value = "2024-07-28"
$ws.Range("C20").Value = "广州·运动番only7.0"
$ws.Range("D20").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws.Range("E20").Value = "2024.07.28 10:00-07.28 17:00"
$ws.Range("F20").Value = 844
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=88473"
$ws.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202407/LtWbGyte1719896481808.jpeg"
$ws.Range("C21").Value = "广州·Le plaisir 第五人格&明日方舟主题同人派对"
$ws.Range("D21").Value = "太和岗路18号负一层 8+1 live house"
$ws.Range("E21").Value = "2024.08.03 13:00-08.03 20:00"
$ws.Range("F21").Value = 76
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=88654"
$ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202407/B3gUl2Gn1720073290274.jpeg"
$ws.Range("C22").Value = "广州·pokemon only PMO 2024-得闲饮茶"
$ws.Range("D22").Value = "较场西路 地王广场"
$ws.Range("E22").Value = "2024.08.03 10:00-08.03 18:00"
$ws.Range("F22").Value = 114
$ws.Range("G22").Value = 69
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=87959"
$ws.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202406/a74dOdcJ1718179975235.jpeg"
$ws.Range("C23").Value = "广州·代号鸢only"
$ws.Range("D23").Value = "会江路 巨大产业园智慧港"
$ws.Range("F23").Value = 289
$ws.Range("G23").Value = 55
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=88224"
$ws.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202406/qBntv1WK1719481529863.jpeg"
$ws.Range("C24").Value = "广州·漫潮动漫游戏嘉年华"
$ws.Range("D24").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws.Range("E24").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F24").Value = 1317
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=86483"
$ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202406/O7TWwoKh1718880739586.jpeg"
$ws.Range("F25").Value = 576
$ws.Range("F26").Value = 443
$ws.Range("F27").Value = 653
$ws.Range("F28").Value = 699
$ws.Range("F30").Value = 54
$ws.Range("F31").Value = 461
$ws.Range("F33").Value = 18
$ws.Range("F34").Value = 495
$ws.Range("F35").Value = 36
$ws.Range("F36").Value = 209
$ws.Range("F37").Value = 905
$ws.Range("F39").Value = 423
$ws.Range("F40").Value = 74
$ws.Range("F41").Value = 133
$ws.Range("F42").Value = 339
$ws.Range("F43").Value = 275
$ws.Range("F44").Value = 1192
$ws.Range("F45").Value = 360
$ws.Range("F48").Value = 14
$ws.Range("F49").Value = 26
